$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 68, 129, 139: session status flips Pending -> Not Recorded (pink fill, same font)
$notRecordedColor = 12695295  # RGB(255,182,193) FFB6C1 in BGR, matches existing "Not Recorded" rows
$ws.Range("A68:I68").Interior.Color = $notRecordedColor
$ws.Range("I68").Value = "Not Recorded"
$ws.Range("A129:I129").Interior.Color = $notRecordedColor
$ws.Range("I129").Value = "Not Recorded"
$ws.Range("A139:I139").Interior.Color = $notRecordedColor
$ws.Range("I139").Value = "Not Recorded"

# Remaining cell-level value updates (recorded-by lists reordered, counts/percentages refreshed)
$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("L7").Value = 14
$ws.Range("L8").Value = 126
$ws.Range("L10").Value = "31.6%"
$ws.Range("S17").Value = "39.6%"
$ws.Range("G18").Value = "shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 16
$ws.Range("G19").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("S19").Value = "33.6%"
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 15
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 16
$ws.Range("S21").Value = "30.8%"
$ws.Range("G24").Value = "rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G40").Value = "shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G41").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("H46").Value = "77/220"
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G54").Value = "basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G58").Value = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G62").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G76").Value = "basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G80").Value = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G84").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G96").Value = "Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("H96").Value = "22/154"
$ws.Range("G98").Value = "basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G118").Value = "Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G120").Value = "basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G134").Value = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H140").Value = "9/224"
$ws.Range("G142").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G150").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("G156").Value = "Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G164").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G172").Value = "wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
